# Add a "PA" column to the codes sheet.
# Inserts a new column before the existing "TB" column (old column F),
# shifting TB/Contact/Pitch one column to the right, then fills the new
# "PA" column: header in F1, and "PA" for every data row that represents
# a completed plate appearance (rows 3-33). Row 2 ("Incomplete Plate
# Appearance") and rows 34-49 ("Advance" type rows) are left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at F; existing F:H (TB, Contact, Pitch) shift to G:I.
$ws.Columns.Item(6).Insert()

# Header for the new column.
$ws.Cells.Item(1, 6).Value = "PA"

# Mark every completed-plate-appearance row (3 through 33) with "PA".
for ($r = 3; $r -le 33; $r++) {
    $ws.Cells.Item($r, 6).Value = "PA"
}
